# Updates the "Jogos da Semana FlashScore" worksheet:
#  - Row 2 (an existing match) is overwritten with new data (a different
#    fixture: Sydney FC vs WS Wanderers) and new odds.
#  - Two brand-new match rows (3 and 4) are appended below it.
# The worksheet's used range grows from A1:BD2 to A1:BD4 as a result.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Row 2 - update existing match (now Sydney FC vs WS Wanderers) ----
$ws.Range("A2").Value = "neJGiEie"
$ws.Range("B2").Value = "23/11/2024"
$ws.Range("C2").Value = "05:35"
$ws.Range("D2").Value = "AUSTRALIA - A-LEAGUE"
$ws.Range("E2").Value = "Sydney FC"
$ws.Range("F2").Value = "WS Wanderers"
$ws.Range("G2").Value = 1.65
$ws.Range("H2").Value = 4.5
$ws.Range("I2").Value = 4.5
$ws.Range("J2").Value = 2.1
$ws.Range("K2").Value = 2.6
$ws.Range("L2").Value = 4.5
$ws.Range("M2").Value = 1.02
$ws.Range("N2").Value = 21
$ws.Range("O2").Value = 1.14
$ws.Range("P2").Value = 5.5
$ws.Range("Q2").Value = 1.44
$ws.Range("R2").Value = 2.75
$ws.Range("S2").Value = 1.25
$ws.Range("T2").Value = 3.75
$ws.Range("U2").Value = 1.5
$ws.Range("V2").Value = 2.5
$ws.Range("W2").Value = 11
$ws.Range("X2").Value = 11
$ws.Range("Y2").Value = 8.5
$ws.Range("Z2").Value = 15
$ws.Range("AA2").Value = 12
$ws.Range("AB2").Value = 17
$ws.Range("AC2").Value = 21
$ws.Range("AD2").Value = 9
$ws.Range("AE2").Value = 13
$ws.Range("AF2").Value = 34
$ws.Range("AG2").Value = 101
$ws.Range("AH2").Value = 19
$ws.Range("AI2").Value = 29
$ws.Range("AJ2").Value = 15
$ws.Range("AK2").Value = 51
$ws.Range("AL2").Value = 29
$ws.Range("AM2").Value = 29
$ws.Range("AN2").Value = 4
$ws.Range("AO2").Value = 8
$ws.Range("AP2").Value = 15
$ws.Range("AQ2").Value = 23
$ws.Range("AR2").Value = 41
$ws.Range("AS2").Value = 81
$ws.Range("AT2").Value = 3.75
$ws.Range("AU2").Value = 7
$ws.Range("AV2").Value = 41
$ws.Range("AW2").Value = 301
$ws.Range("AX2").Value = 7
$ws.Range("AY2").Value = 21
$ws.Range("AZ2").Value = 23
$ws.Range("BA2").Value = 67
$ws.Range("BB2").Value = 67
$ws.Range("BC2").Value = 126
$ws.Range("BD2").Value = 151

# ---- Row 3 - new match (Barito Putera vs Persita) ----
$ws.Range("A3").Value = "tSFajqxe"
$ws.Range("B3").Value = "23/11/2024"
$ws.Range("C3").Value = "05:30"
$ws.Range("D3").Value = "INDONESIA - LIGA 1"
$ws.Range("E3").Value = "Barito Putera"
$ws.Range("F3").Value = "Persita"
$ws.Range("G3").Value = 1.88
$ws.Range("H3").Value = 3.2
$ws.Range("I3").Value = 4.1
$ws.Range("J3").Value = 2.45
$ws.Range("K3").Value = 2.07
$ws.Range("L3").Value = 4.35
$ws.Range("M3").Value = 1.01
$ws.Range("N3").Value = 7.9
$ws.Range("O3").Value = 1.34
$ws.Range("P3").Value = 2.75
$ws.Range("Q3").Value = 2
$ws.Range("R3").Value = 1.65
$ws.Range("S3").Value = 1.39
$ws.Range("T3").Value = 2.57
$ws.Range("U3").Value = 1.83
$ws.Range("V3").Value = 1.78
$ws.Range("W3").Value = 6.2
$ws.Range("X3").Value = 8.25
$ws.Range("Y3").Value = 8.5
$ws.Range("Z3").Value = 16
$ws.Range("AA3").Value = 16.5
$ws.Range("AB3").Value = 30
$ws.Range("AC3").Value = 8.25
$ws.Range("AD3").Value = 6.2
$ws.Range("AE3").Value = 15.5
$ws.Range("AF3").Value = 80
$ws.Range("AG3").Value = 700
$ws.Range("AH3").Value = 10.75
$ws.Range("AI3").Value = 23
$ws.Range("AJ3").Value = 13.5
$ws.Range("AK3").Value = 70
$ws.Range("AL3").Value = 40
$ws.Range("AM3").Value = 45
$ws.Range("AN3").Value = 3.7
$ws.Range("AO3").Value = 9.5
$ws.Range("AP3").Value = 18
$ws.Range("AQ3").Value = 35
$ws.Range("AR3").Value = 65
$ws.Range("AS3").Value = 250
$ws.Range("AT3").Value = 2.52
$ws.Range("AU3").Value = 6.9
$ws.Range("AV3").Value = 60
# AW3 is blank in the source data
$ws.Range("AX3").Value = 5.8
$ws.Range("AY3").Value = 22
$ws.Range("AZ3").Value = 26
$ws.Range("BA3").Value = 120
$ws.Range("BB3").Value = 150
$ws.Range("BC3").Value = 300
# BD3 is blank in the source data

# ---- Row 4 - new match (Persik Kediri vs PSIS Semarang) ----
$ws.Range("A4").Value = "Y1Cih57r"
$ws.Range("B4").Value = "23/11/2024"
$ws.Range("C4").Value = "05:30"
$ws.Range("D4").Value = "INDONESIA - LIGA 1"
$ws.Range("E4").Value = "Persik Kediri"
$ws.Range("F4").Value = "PSIS Semarang"
$ws.Range("G4").Value = 1.7
$ws.Range("H4").Value = 3.35
$ws.Range("I4").Value = 5
$ws.Range("J4").Value = 2.25
$ws.Range("K4").Value = 2.12
$ws.Range("L4").Value = 5
$ws.Range("M4").Value = 1.03
$ws.Range("N4").Value = 6.75
$ws.Range("O4").Value = 1.32
$ws.Range("P4").Value = 2.82
$ws.Range("Q4").Value = 1.98
$ws.Range("R4").Value = 1.65
$ws.Range("S4").Value = 1.4
$ws.Range("T4").Value = 2.52
$ws.Range("U4").Value = 1.87
$ws.Range("V4").Value = 1.75
$ws.Range("W4").Value = 6
$ws.Range("X4").Value = 7.3
$ws.Range("Y4").Value = 8.25
$ws.Range("Z4").Value = 13
$ws.Range("AA4").Value = 14.5
$ws.Range("AB4").Value = 30
$ws.Range("AC4").Value = 8.75
$ws.Range("AD4").Value = 6.6
$ws.Range("AE4").Value = 16.5
$ws.Range("AF4").Value = 80
$ws.Range("AG4").Value = 700
$ws.Range("AH4").Value = 12.5
$ws.Range("AI4").Value = 29
$ws.Range("AJ4").Value = 16
$ws.Range("AK4").Value = 100
$ws.Range("AL4").Value = 50
$ws.Range("AM4").Value = 55
$ws.Range("AN4").Value = 3.45
$ws.Range("AO4").Value = 8.25
$ws.Range("AP4").Value = 17.5
$ws.Range("AQ4").Value = 28
$ws.Range("AR4").Value = 60
$ws.Range("AS4").Value = 250
$ws.Range("AT4").Value = 2.55
$ws.Range("AU4").Value = 7.2
$ws.Range("AV4").Value = 65
# AW4 is blank in the source data
$ws.Range("AX4").Value = 6.5
$ws.Range("AY4").Value = 28
$ws.Range("AZ4").Value = 30
$ws.Range("BA4").Value = 175
$ws.Range("BB4").Value = 175
$ws.Range("BC4").Value = 400
# BD4 is blank in the source data
